$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped cryptos data.
# D-column values that look like plain decimals (e.g. "213.20") are prefixed
# with a leading apostrophe so Excel stores them as literal text (preserving
# trailing zeros) instead of silently converting them to numbers; values that
# already contain two dots (e.g. "26.454.65") are never auto-parsed as
# numbers, so they're set as-is.
$ws.Range("D2").Value = "26.454.65"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.626.94"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'213.20"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.0624"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'18.91"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.853.60"
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("D13").Value = "1.624.23"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'64.80"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "26.526.17"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'215.07"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'6.27"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'2.00"
$ws.Range("E24").Value = "  +4.25%  "
$ws.Range("D25").Value = "'148.49"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").Value = "'2.38"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "1.218.93"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").Value = "'0.506"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("D42").Value = "'0.793"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'5.36"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.764.54"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "'93.10"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").Value = "'54.85"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  +0.14%  "
